$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.478.02"
$ws.Range("E2").Value = "  +1.15%  "

$ws.Range("D3").Value = "1.907.38"
$ws.Range("E3").Value = "  +2.82%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.43"
$ws.Range("E5").Value = "  +3.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.636"
$ws.Range("E6").Value = "  +2.24%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.06"
$ws.Range("E8").Value = "  -1.34%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.339"
$ws.Range("E9").Value = "  +2.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0704"
$ws.Range("E10").Value = "  +1.47%  "

$ws.Range("E11").Value = "  +0.77%  "

$ws.Range("D12").Value = "2.182.98"
$ws.Range("E12").Value = "  +2.73%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "12.32"
$ws.Range("E13").Value = "  +8.14%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.934.67"
$ws.Range("E14").Value = "  +4.14%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.695"
$ws.Range("E15").Value = "  +2.34%  "

$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.81"
$ws.Range("E16").Value = "  +2.73%  "

$ws.Range("D17").Value = "35.540.07"
$ws.Range("E17").Value = "  +1.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "71.90"
$ws.Range("E18").Value = "  +2.22%  "

$ws.Range("E19").Value = "  +3.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "243.40"
$ws.Range("E20").Value = "  +0.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.70"
$ws.Range("E21").Value = "  +4.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.84"
$ws.Range("E22").Value = "  +1.79%  "

$ws.Range("E23").Value = "  +0.21%  "

$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.14"
$ws.Range("E25").Value = "  +0.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.20"
$ws.Range("E26").Value = "  +16.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.59"
$ws.Range("E27").Value = "  +8.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.00"
$ws.Range("E28").Value = "  +2.06%  "

$ws.Range("E29").Value = "  +1.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.978"
$ws.Range("E30").Value = "  +25.74%  "

$ws.Range("E31").Value = "  +2.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.11"
$ws.Range("E32").Value = "  +3.00%  "

$ws.Range("E33").Value = "  +0.21%  "

$ws.Range("E34").Value = "  +4.44%  "

$ws.Range("E35").Value = "  +7.76%  "

$ws.Range("E36").Value = "  +0.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.34"
$ws.Range("E37").Value = "  +7.99%  "

$ws.Range("E38").Value = "  +2.90%  "

$ws.Range("E39").Value = "  +15.98%  "

$ws.Range("E40").Value = "  +1.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "91.54"
$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.68"
$ws.Range("E42").Value = "  +4.55%  "

$ws.Range("D43").Value = "1.350.44"
$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "50.26"
$ws.Range("E44").Value = "  +45.04%  "

$ws.Range("E45").Value = "  +3.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.81"
$ws.Range("E46").Value = "  +0.21%  "

$ws.Range("E47").Value = "  +0.15%  "

$ws.Range("E48").Value = "  -0.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.65"
$ws.Range("E49").Value = "  +3.58%  "

$ws.Range("D50").Value = "2.094.11"
$ws.Range("E50").Value = "  +2.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0693"
$ws.Range("E51").Value = "  +1.85%  "
